# Updated log and error mechanism
# Target Column values on the T_EMP sheet are updated, and the active
# selection moves from G2 to G3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("T_EMP")

$ws.Range("G2").Value = "T_EMP_ID1212"
$ws.Range("G3").Value = "EMP_ID345"

$ws.Range("G3").Select()
